# Generate Report for Archive
# Updates the localization status workbook: the files
# "2d9b2211-ee9e-441d-b55c-e43b682ff33b.md" and
# "de797916-d049-4ef7-8a19-e4adb1921ace.md" move from status
# "Ready for handoff" to "In Translation" on the Overview sheet
# (zh-cn / de-de columns) as well as on the per-locale zh-cn and
# de-de report sheets (Status column).

$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")

# Row 3 -> 2d9b2211-ee9e-441d-b55c-e43b682ff33b.md
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Row 4 -> de797916-d049-4ef7-8a19-e4adb1921ace.md
$wsOverview.Range("E4").Value = "In Translation"
$wsOverview.Range("F4").Value = "In Translation"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")

# Row 3 -> 2d9b2211-ee9e-441d-b55c-e43b682ff33b.md
$wsZhCn.Range("C3").Value = "In Translation"
# Row 4 -> de797916-d049-4ef7-8a19-e4adb1921ace.md
$wsZhCn.Range("C4").Value = "In Translation"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")

# Row 3 -> 2d9b2211-ee9e-441d-b55c-e43b682ff33b.md
$wsDeDe.Range("C3").Value = "In Translation"
# Row 4 -> de797916-d049-4ef7-8a19-e4adb1921ace.md
$wsDeDe.Range("C4").Value = "In Translation"
